$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-11-02 Thursday" "2023-11-03 Friday"

Replace-Text "41×34=" "23×38="
Replace-Text "39×67=" "51×40="
Replace-Text "48×72=" "91×11="
Replace-Text "30×44=" "18×56="
Replace-Text "90×41=" "51×42="

Replace-Text "30×59=" "32×88="
Replace-Text "64×26=" "51×84="
Replace-Text "16×84=" "95×56="
Replace-Text "54×67=" "43×57="
Replace-Text "79×31=" "77×15="

Replace-Text "25×40=" "17×20="
Replace-Text "32×98=" "91×62="
Replace-Text "13×72=" "56×23="
Replace-Text "93×32=" "38×15="
Replace-Text "42×88=" "59×49="

Replace-Text "14×35=" "46×54="
Replace-Text "58×43=" "57×85="
Replace-Text "84×38=" "80×35="
Replace-Text "48×56=" "94×81="
Replace-Text "54×41=" "77×29="

Replace-Text "52×89=" "41×19="
Replace-Text "12×74=" "43×49="
Replace-Text "64×67=" "48×12="
Replace-Text "77×61=" "40×75="
Replace-Text "50×68=" "82×44="
